$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 67: phone number in A67 was stored as text; normalize it to a real number
$ws.Cells.Item(67, 1).Value = 71717170

# Row 68: append new payment record -> 71717170 (Cash) 2025-08-20T08:18:12
# Phone number is kept as text (as in the source data), so force text formatting
# before assigning the numeric-looking string, then drop the residual format
# so no extra style is attached to the cell.
$ws.Cells.Item(68, 1).NumberFormat = "@"
$ws.Cells.Item(68, 1).Value = "71717170"
$ws.Cells.Item(68, 1).ClearFormats()

$ws.Cells.Item(68, 3).Value = "Cash"
$ws.Cells.Item(68, 4).Value = "2025-08-20T08:18:12"
$ws.Cells.Item(68, 5).Value = 137
$ws.Cells.Item(68, 7).Value = 116.45
$ws.Cells.Item(68, 8).Value = 20.55
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
